# expenditure_of_time.xlsx — "Done akutaliesiert + offene Punkte ergänzt"
#
# Row 14 ("time costs" sheet):
#  - G14 text gets two extra lines appended (login fix + duell-request popup
#    workaround notes).
#  - H14 text drops the now-obsolete "login führt..." line and gets a
#    trailing "Präsentation" line appended.
#  - A brand-new I14 cell is added with "Zu Beachten / Zu Besprechen" notes.
#  - Row 14 grows taller (255 -> 285) to fit the extra text.
#  - The remembered selection moves from G15 to G14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("time costs")

$g14 = "Graue Vierecke anstatt `"Verdeckt`"`nSpielaufgabe disabled Spielen button bei gegner`nApp speichert random categories pro Runde für sich, um nicht mit Zurück-Vor-Navigation wieder neue Kategorien zu bekommen.`nServer liefert 5 letzte Spiele (beendet und aufgegeben)`nLogin functioniert nun (führte manchmal nicht zum home screen)`nWorkaround für Duellanfragen werden u.U mehrmals im Hauptmenü angezeigt! (popUp) wurde implementiert"

$h14 = "Buttons nur auslösen, wenn Anfang & ende des toches drauf sind!`nRandomEnemy (serverseitig!)`n`nBei SpielEnde Benachrichtigung & Ändern des SpielenButtons`nLogo anzeigen (inapp & icon!)`nFrage Buttons schrift zu klein nach Auswertung!`nBeendete Spiele (letzte 5) in Sync mit liefern & in Hauptmenü anzeigen.`nDuellanfragen werden u.U mehrmals im Hauptmenü angezeigt! (popUp)`ntastaturinput-enter --> Aktion auf screen! (login/suchen..)`nweiterbutton durch swipe ersetzen`naktualisieren buttons in RÜ & home in navigationbar`nPräsentation"

$i14 = "Zu Beachten: `n1. Anzahl der Serveraufrufe beschränken, wichtig für späteren verlauf`n(Mögliche Überlastung des Servers)`n2. Einseitig Funktion von Steroids vermeiden (Bild in Titelzeile ist nicht für Android verfügbar) `nZu Besprechen: weitere Aufgabenverteilung, was wollen wir noch implementieren/was brauchen wir noch für den Prototyp`n"

$ws.Range("G14").Value = $g14
$ws.Range("I14").Value = $i14
$ws.Range("H14").Value = $h14

$ws.Range("I14").WrapText = $true
$ws.Range("I14").HorizontalAlignment = -4131
$ws.Range("I14").VerticalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 285

$ws.Range("G14").Select()
